$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Stock count decrements on existing rows ---
$ws.Range("B14").Value = 4
$ws.Range("B107").Value = 6
$ws.Range("B220").Value = 0
$ws.Range("B247").Value = 1
$ws.Range("B256").Value = 0
$ws.Range("B384").Value = 1

# --- Append new rows (419-452) for the new "Bricks" minifigure category ---
# Copy formatting from the last existing row (418) down for each new row first,
# which mirrors how Excel extends formatting for appended rows.
$ws.Rows(418).Copy()
$ws.Rows(419).Insert(-4121)
$ws.Range("A419").Value = "Pikachu MINI Bricks"
$ws.Range("B419").Value = 1
$ws.Range("C419").Value = 200
$ws.Range("D419").Value = "pikachu bricks.jpg"
$ws.Range("E419").Value = "Bricks"

$ws.Rows(419).Copy()
$ws.Rows(420).Insert(-4121)
$ws.Range("A420").Value = "Hello Kitty MINI Bricks"
$ws.Range("B420").Value = 1
$ws.Range("C420").Value = 200
$ws.Range("D420").Value = "Kitty bricks.jpg"
$ws.Range("E420").Value = "Bricks"

$ws.Rows(420).Copy()
$ws.Rows(421).Insert(-4121)
$ws.Range("A421").Value = "Kuromi MINI Bricks"
$ws.Range("B421").Value = 1
$ws.Range("C421").Value = 200
$ws.Range("D421").Value = "kuromi bricks.jpg"
$ws.Range("E421").Value = "Bricks"

$ws.Rows(421).Copy()
$ws.Rows(422).Insert(-4121)
$ws.Range("A422").Value = "Jigglypuff MINI Bricks"
$ws.Range("B422").Value = 1
$ws.Range("C422").Value = 200
$ws.Range("D422").Value = "jigglypuff bricks.jpg"
$ws.Range("E422").Value = "Bricks"

$ws.Rows(422).Copy()
$ws.Rows(423).Insert(-4121)
$ws.Range("A423").Value = "Charizard MINI Bricks"
$ws.Range("B423").Value = 1
$ws.Range("C423").Value = 200
$ws.Range("D423").Value = "charizard bricks.jpg"
$ws.Range("E423").Value = "Bricks"

$ws.Rows(423).Copy()
$ws.Rows(424).Insert(-4121)
$ws.Range("A424").Value = "Snorlax MINI Bricks"
$ws.Range("B424").Value = 1
$ws.Range("C424").Value = 200
$ws.Range("D424").Value = "snorlax bricks.jpg"
$ws.Range("E424").Value = "Bricks"

$ws.Rows(424).Copy()
$ws.Rows(425).Insert(-4121)
$ws.Range("A425").Value = "Eevee MINI Bricks"
$ws.Range("B425").Value = 1
$ws.Range("C425").Value = 200
$ws.Range("D425").Value = "eevee bricks.jpg"
$ws.Range("E425").Value = "Bricks"

$ws.Rows(425).Copy()
$ws.Rows(426).Insert(-4121)
$ws.Range("A426").Value = "Venasaur MINI Bricks"
$ws.Range("B426").Value = 1
$ws.Range("C426").Value = 200
$ws.Range("D426").Value = "venasaur bricks.jpg"
$ws.Range("E426").Value = "Bricks"

$ws.Rows(426).Copy()
$ws.Rows(427).Insert(-4121)
$ws.Range("A427").Value = "Luffy MINI Bricks"
$ws.Range("B427").Value = 2
$ws.Range("C427").Value = 200
$ws.Range("D427").Value = "luffy bricks.jpg"
$ws.Range("E427").Value = "Bricks"

$ws.Rows(427).Copy()
$ws.Rows(428).Insert(-4121)
$ws.Range("A428").Value = "Zoro MINI Bricks"
$ws.Range("B428").Value = 2
$ws.Range("C428").Value = 200
$ws.Range("D428").Value = "zoro bricks.jpg"
$ws.Range("E428").Value = "Bricks"

$ws.Rows(428).Copy()
$ws.Rows(429).Insert(-4121)
$ws.Range("A429").Value = "Sanji MINI Bricks"
$ws.Range("B429").Value = 2
$ws.Range("C429").Value = 200
$ws.Range("D429").Value = "sanji bricks.jpg"
$ws.Range("E429").Value = "Bricks"

$ws.Rows(429).Copy()
$ws.Rows(430).Insert(-4121)
$ws.Range("A430").Value = "Chopper MINI Bricks"
$ws.Range("B430").Value = 1
$ws.Range("C430").Value = 200
$ws.Range("D430").Value = "chopper bricks.jpg"
$ws.Range("E430").Value = "Bricks"

$ws.Rows(430).Copy()
$ws.Rows(431).Insert(-4121)
$ws.Range("A431").Value = "Robin MINI Bricks"
$ws.Range("B431").Value = 1
$ws.Range("C431").Value = 200
$ws.Range("D431").Value = "robin bricks.jpg"
$ws.Range("E431").Value = "Bricks"

$ws.Rows(431).Copy()
$ws.Rows(432).Insert(-4121)
$ws.Range("A432").Value = "Ace MINI Bricks"
$ws.Range("B432").Value = 1
$ws.Range("C432").Value = 200
$ws.Range("D432").Value = "ace bricks.jpg"
$ws.Range("E432").Value = "Bricks"

$ws.Rows(432).Copy()
$ws.Rows(433).Insert(-4121)
$ws.Range("A433").Value = "Naruto MINI Bricks"
$ws.Range("B433").Value = 2
$ws.Range("C433").Value = 200
$ws.Range("D433").Value = "naruto bricks.jpg"
$ws.Range("E433").Value = "Bricks"

$ws.Rows(433).Copy()
$ws.Rows(434).Insert(-4121)
$ws.Range("A434").Value = "Sasuke MINI Bricks"
$ws.Range("B434").Value = 2
$ws.Range("C434").Value = 200
$ws.Range("D434").Value = "sasuke bricks.jpg"
$ws.Range("E434").Value = "Bricks"

$ws.Rows(434).Copy()
$ws.Rows(435).Insert(-4121)
$ws.Range("A435").Value = "Itachi MINI Bricks"
$ws.Range("B435").Value = 2
$ws.Range("C435").Value = 200
$ws.Range("D435").Value = "itachi bricks.jpg"
$ws.Range("E435").Value = "Bricks"

$ws.Rows(435).Copy()
$ws.Rows(436).Insert(-4121)
$ws.Range("A436").Value = "Kakashi MINI Bricks"
$ws.Range("B436").Value = 2
$ws.Range("C436").Value = 200
$ws.Range("D436").Value = "kakashi bricks.jpg"
$ws.Range("E436").Value = "Bricks"

$ws.Rows(436).Copy()
$ws.Rows(437).Insert(-4121)
$ws.Range("A437").Value = "Goku MINI Bricks"
$ws.Range("B437").Value = 2
$ws.Range("C437").Value = 200
$ws.Range("D437").Value = "goku bricks.jpg"
$ws.Range("E437").Value = "Bricks"

$ws.Rows(437).Copy()
$ws.Rows(438).Insert(-4121)
$ws.Range("A438").Value = "Tanjiro MINI Bricks"
$ws.Range("B438").Value = 1
$ws.Range("C438").Value = 200
$ws.Range("D438").Value = "tanjiro bricks.jpg"
$ws.Range("E438").Value = "Bricks"

$ws.Rows(438).Copy()
$ws.Rows(439).Insert(-4121)
$ws.Range("A439").Value = "Nezuko MINI Bricks"
$ws.Range("B439").Value = 2
$ws.Range("C439").Value = 200
$ws.Range("D439").Value = "nezuko bricks.jpg"
$ws.Range("E439").Value = "Bricks"

$ws.Rows(439).Copy()
$ws.Rows(440).Insert(-4121)
$ws.Range("A440").Value = "Zenitsu MINI Bricks"
$ws.Range("B440").Value = 2
$ws.Range("C440").Value = 200
$ws.Range("D440").Value = "zenitsu bricks.jpg"
$ws.Range("E440").Value = "Bricks"

$ws.Rows(440).Copy()
$ws.Rows(441).Insert(-4121)
$ws.Range("A441").Value = "Giyuu MINI Bricks"
$ws.Range("B441").Value = 1
$ws.Range("C441").Value = 200
$ws.Range("D441").Value = "giyuu bricks.jpg"
$ws.Range("E441").Value = "Bricks"

$ws.Rows(441).Copy()
$ws.Rows(442).Insert(-4121)
$ws.Range("A442").Value = "Rengoku MINI Bricks"
$ws.Range("B442").Value = 1
$ws.Range("C442").Value = 200
$ws.Range("D442").Value = "rengoku bricks.jpg"
$ws.Range("E442").Value = "Bricks"

$ws.Rows(442).Copy()
$ws.Rows(443).Insert(-4121)
$ws.Range("A443").Value = "Tengen MINI Bricks"
$ws.Range("B443").Value = 1
$ws.Range("C443").Value = 200
$ws.Range("D443").Value = "tengen bricks.jpg"
$ws.Range("E443").Value = "Bricks"

$ws.Rows(443).Copy()
$ws.Rows(444).Insert(-4121)
$ws.Range("A444").Value = "Shinobou MINI Bricks"
$ws.Range("B444").Value = 1
$ws.Range("C444").Value = 200
$ws.Range("D444").Value = "shinobou bricks.jpg"
$ws.Range("E444").Value = "Bricks"

$ws.Rows(444).Copy()
$ws.Rows(445).Insert(-4121)
$ws.Range("A445").Value = "Muichiro MINI Bricks"
$ws.Range("B445").Value = 1
$ws.Range("C445").Value = 200
$ws.Range("D445").Value = "muichiro bricks.jpg"
$ws.Range("E445").Value = "Bricks"

$ws.Rows(445).Copy()
$ws.Rows(446).Insert(-4121)
$ws.Range("A446").Value = "Ash Ketchum [Hoenn]"
$ws.Range("B446").Value = 1
$ws.Range("C446").Value = 200
$ws.Range("D446").Value = "ash hoenn.jpg"
$ws.Range("E446").Value = "Pokemon"

$ws.Rows(446).Copy()
$ws.Rows(447).Insert(-4121)
$ws.Range("A447").Value = "Ash Ketchum [Sinnoh]"
$ws.Range("B447").Value = 1
$ws.Range("C447").Value = 200
$ws.Range("D447").Value = "ash sinnoh.jpg"
$ws.Range("E447").Value = "Pokemon"

$ws.Rows(447).Copy()
$ws.Rows(448).Insert(-4121)
$ws.Range("A448").Value = "Ash Ketchum [Unova]"
$ws.Range("B448").Value = 3
$ws.Range("C448").Value = 200
$ws.Range("D448").Value = "ash unova.jpg"
$ws.Range("E448").Value = "Pokemon"

$ws.Rows(448).Copy()
$ws.Rows(449).Insert(-4121)
$ws.Range("A449").Value = "Ash Ketchum [Kalos]"
$ws.Range("B449").Value = 3
$ws.Range("C449").Value = 200
$ws.Range("D449").Value = "ash kalos.jpg"
$ws.Range("E449").Value = "Pokemon"

$ws.Rows(449).Copy()
$ws.Rows(450).Insert(-4121)
$ws.Range("A450").Value = "Misty"
$ws.Range("B450").Value = 1
$ws.Range("C450").Value = 150
$ws.Range("D450").Value = "misty.jpg"
$ws.Range("E450").Value = "Pokemon"

$ws.Rows(450).Copy()
$ws.Rows(451).Insert(-4121)
$ws.Range("A451").Value = "Serena"
$ws.Range("B451").Value = 4
$ws.Range("C451").Value = 150
$ws.Range("D451").Value = "serena.jpg"
$ws.Range("E451").Value = "Pokemon"

$ws.Rows(451).Copy()
$ws.Rows(452).Insert(-4121)
$ws.Range("A452").Value = "Calem"
$ws.Range("B452").Value = 2
$ws.Range("C452").Value = 150
$ws.Range("D452").Value = "calem.jpg"
$ws.Range("E452").Value = "Pokemon"

# --- Update selection / scroll position to mirror final editing state ---
$ws.Range("B384").Select()
$excel.ActiveWindow.ScrollRow = 367
$excel.ActiveWindow.ScrollColumn = 1